# order extensibility and receipt ui
# Fill in row 12 of the menu list with a new "fresh burger" item.
# Columns: A=id, B=name, C=price, D=branch, E=category, F=description

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 12

# id (uuid) - plain text
$cA = $ws.Cells.Item($row, 1)
$cA.Value = "04e97c9d-4c19-4a1a-9046-b975f194ac6a"
$cA.Style = "Normal"

# name - reuse existing text "burger"
$cB = $ws.Cells.Item($row, 2)
$cB.Value = "burger"
$cB.Style = "Normal"

# price - force as text (not an auto-converted number) using a leading
# quote prefix, then clear the resulting style so the cell keeps default
# formatting while remaining a text value.
$cC = $ws.Cells.Item($row, 3)
$cC.Formula = "'4.5"
$cC.Style = "Normal"

# branch - reuse existing text "NTU"
$cD = $ws.Cells.Item($row, 4)
$cD.Value = "NTU"
$cD.Style = "Normal"

# category - reuse existing text "Test"
$cE = $ws.Cells.Item($row, 5)
$cE.Value = "Test"
$cE.Style = "Normal"

# description - new text
$cF = $ws.Cells.Item($row, 6)
$cF.Value = "fresh burger"
$cF.Style = "Normal"
